# Update cached numeric results (currentAveragePrice / Leve price & profit
# columns) across all 8 job sheets per the scheduled market-data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 156.3
$ws.Range("I6").Value = 148
$ws.Range("K6").Value = 444
$ws.Range("M6").Value = -332
$ws.Range("H129").Value = 4214.3335
$ws.Range("J129").Value = 8947.6
$ws.Range("L129").Value = 26842.8
$ws.Range("N129").Value = -36842.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 1000000000
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 1000000000
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 1000000000
$ws.Range("N6").Value = -1000000346
$ws.Range("M6").ClearContents()
$ws.Range("H32").Value = 7509.8125
$ws.Range("I32").Value = 7509.8125
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 7509.8125
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -7222.8125
$ws.Range("N32").ClearContents()
$ws.Range("H74").Value = 1847.1818
$ws.Range("I74").Value = 1691.4736
$ws.Range("K74").Value = 1691.4736
$ws.Range("M74").Value = -817.4736
$ws.Range("H77").Value = 1847.1818
$ws.Range("I77").Value = 1691.4736
$ws.Range("K77").Value = 8457.368
$ws.Range("M77").Value = -4089.368

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3282.162
$ws.Range("I20").Value = 2679.0908
$ws.Range("K20").Value = 2679.0908
$ws.Range("M20").Value = -2432.0908
$ws.Range("H86").Value = 6498.6
$ws.Range("I86").Value = 1352.5
$ws.Range("J86").Value = 9929.333
$ws.Range("K86").Value = 1352.5
$ws.Range("L86").Value = 9929.333
$ws.Range("M86").Value = -229.5
$ws.Range("N86").Value = -12175.333
$ws.Range("H89").Value = 6498.6
$ws.Range("I89").Value = 1352.5
$ws.Range("J89").Value = 9929.333
$ws.Range("K89").Value = 6762.5
$ws.Range("L89").Value = 49646.665
$ws.Range("M89").Value = -1146.5
$ws.Range("N89").Value = -60878.665
$ws.Range("H122").Value = 46208.1
$ws.Range("J122").Value = 46000
$ws.Range("L122").Value = 46000
$ws.Range("N122").Value = -55800
$ws.Range("H134").Value = 33334170
$ws.Range("I134").Value = 1256
$ws.Range("K134").Value = 3768
$ws.Range("M134").Value = -1233

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 34486170
$ws.Range("I31").Value = 71431190
$ws.Range("J31").Value = 4151
$ws.Range("K31").Value = 71431190
$ws.Range("L31").Value = 4151
$ws.Range("M31").Value = -71430895
$ws.Range("N31").Value = -4741
$ws.Range("H34").Value = 34486170
$ws.Range("I34").Value = 71431190
$ws.Range("J34").Value = 4151
$ws.Range("K34").Value = 71431190
$ws.Range("L34").Value = 4151
$ws.Range("M34").Value = -71430988
$ws.Range("N34").Value = -4555
$ws.Range("H64").Value = 74782.664
$ws.Range("J64").Value = 74782.664
$ws.Range("L64").Value = 74782.664
$ws.Range("N64").Value = -75278.664
$ws.Range("H67").Value = 74782.664
$ws.Range("J67").Value = 74782.664
$ws.Range("L67").Value = 74782.664
$ws.Range("N67").Value = -76498.664
$ws.Range("H134").Value = 2827.125
$ws.Range("I134").Value = 2603
$ws.Range("J134").Value = 3499.5
$ws.Range("K134").Value = 7809
$ws.Range("L134").Value = 10498.5
$ws.Range("M134").Value = -5274
$ws.Range("N134").Value = -15568.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 6700.364
$ws.Range("J109").Value = 17777.666
$ws.Range("L109").Value = 53332.99800000001
$ws.Range("N109").Value = -55412.99800000001
$ws.Range("H117").Value = 3758.3333
$ws.Range("I117").Value = 314
$ws.Range("K117").Value = 942
$ws.Range("M117").Value = 2500
$ws.Range("H132").Value = 1899.4286
$ws.Range("J132").Value = 1399.6666
$ws.Range("L132").Value = 12596.9994
$ws.Range("N132").Value = -17656.9994
$ws.Range("H133").Value = 21613
$ws.Range("I133").Value = 16696.6
$ws.Range("K133").Value = 50089.8
$ws.Range("M133").Value = -45029.8
$ws.Range("H136").Value = 19166.5
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H139").Value = 5928.0713
$ws.Range("I139").Value = 3167.7646
$ws.Range("K139").Value = 9503.2938
$ws.Range("M139").Value = -4363.293799999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 236.5
$ws.Range("I2").Value = 218.5
$ws.Range("J2").Value = 272.5
$ws.Range("K2").Value = 218.5
$ws.Range("L2").Value = 272.5
$ws.Range("M2").Value = -105.5
$ws.Range("N2").Value = -498.5
$ws.Range("H39").Value = 99999
$ws.Range("J39").Value = 99999
$ws.Range("L39").Value = 99999
$ws.Range("N39").Value = -101063
$ws.Range("H70").Value = 14934.1
$ws.Range("J70").Value = 14922.167
$ws.Range("L70").Value = 14922.167
$ws.Range("N70").Value = -15462.167
$ws.Range("H73").Value = 14934.1
$ws.Range("J73").Value = 14922.167
$ws.Range("L73").Value = 14922.167
$ws.Range("N73").Value = -16794.167
$ws.Range("H132").Value = 3129248.5
$ws.Range("I132").Value = 4365.1665
$ws.Range("K132").Value = 13095.4995
$ws.Range("M132").Value = -10565.4995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7173.154
$ws.Range("I7").Value = 6920.615
$ws.Range("J7").Value = 7678.231
$ws.Range("K7").Value = 6920.615
$ws.Range("L7").Value = 7678.231
$ws.Range("M7").Value = -6808.615
$ws.Range("N7").Value = -7902.231
$ws.Range("H98").Value = 90177
$ws.Range("J98").Value = 90177
$ws.Range("L98").Value = 90177
$ws.Range("N98").Value = -96167
$ws.Range("H126").Value = 7173.154
$ws.Range("I126").Value = 6920.615
$ws.Range("J126").Value = 7678.231
$ws.Range("K126").Value = 20761.845
$ws.Range("L126").Value = 23034.693
$ws.Range("M126").Value = -18291.845
$ws.Range("N126").Value = -27974.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 12966.556
$ws.Range("I62").Value = 4933.3335
$ws.Range("J62").Value = 16983.166
$ws.Range("K62").Value = 4933.3335
$ws.Range("L62").Value = 16983.166
$ws.Range("M62").Value = -4309.3335
$ws.Range("N62").Value = -18231.166
$ws.Range("H65").Value = 12966.556
$ws.Range("I65").Value = 4933.3335
$ws.Range("J65").Value = 16983.166
$ws.Range("K65").Value = 24666.6675
$ws.Range("L65").Value = 84915.83
$ws.Range("M65").Value = -21546.6675
$ws.Range("N65").Value = -91155.83
$ws.Range("H81").Value = 1323.8
$ws.Range("I81").Value = 1217.375
$ws.Range("K81").Value = 2434.75
$ws.Range("M81").Value = -1373.75
$ws.Range("H84").Value = 1323.8
$ws.Range("I84").Value = 1217.375
$ws.Range("K84").Value = 12173.75
$ws.Range("M84").Value = -6869.75
$ws.Range("H105").Value = 39305.5
$ws.Range("J105").Value = 39305.5
$ws.Range("L105").Value = 39305.5
$ws.Range("N105").Value = -46293.5
$ws.Range("H122").Value = 1344.2593
$ws.Range("I122").Value = 1011.38464
$ws.Range("K122").Value = 3034.15392
$ws.Range("M122").Value = -584.1539199999997
$ws.Range("H132").Value = 1113404.4
$ws.Range("I132").Value = 2283.7693
$ws.Range("K132").Value = 6851.3079
$ws.Range("M132").Value = -4321.3079
